$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = "13-Oct"
$ws.Range("B4").Value = "АОР"
$ws.Range("C4").Value = "Посев"
$ws.Range("D4").Value = "Ячмень яровой"
$ws.Range("E4").Value = 85
$ws.Range("F4").Value = 1250
$ws.Range("G4").Value = 420
$ws.Range("H4").Value = 5800

# Row 5
$ws.Range("A5").Value = "14-Oct"
$ws.Range("B5").Value = "Центральное"
$ws.Range("C5").Value = "Уборка"
$ws.Range("D5").Value = "Кукуруза"
$ws.Range("E5").Value = 210
$ws.Range("F5").Value = 3150
$ws.Range("G5").Value = 1250
$ws.Range("H5").Value = 18700

# Row 6
$ws.Range("A6").Value = "12-Oct"
$ws.Range("B6").Value = "АОР"
$ws.Range("C6").Value = "Внесение минеральных удобрений"
$ws.Range("D6").Value = "Пшеница озимая"
$ws.Range("E6").Value = 149
$ws.Range("F6").Value = 7264

# Row 7
$ws.Range("A7").Value = "13-Oct"
$ws.Range("B7").Value = "АОР"
$ws.Range("C7").Value = "Посев"
$ws.Range("D7").Value = "Ячмень яровой"
$ws.Range("E7").Value = 85
$ws.Range("F7").Value = 1250
$ws.Range("G7").Value = 420
$ws.Range("H7").Value = 5800

# Row 8
$ws.Range("A8").Value = "14-Oct"
$ws.Range("B8").Value = "Центральное"
$ws.Range("C8").Value = "Уборка"
$ws.Range("D8").Value = "Кукуруза"
$ws.Range("E8").Value = 210
$ws.Range("F8").Value = 3150
$ws.Range("G8").Value = 1250
$ws.Range("H8").Value = 18700

# Apply the same formatting as row 3 (style s="2") to the new rows 4-8
$ws.Range("A3:H3").Copy()
$ws.Range("A4:H8").PasteSpecial(-4122)
